$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "About"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item('About')
$ws1.UsedRange.Clear()

$ws1.Range('A1').Value2 = 'RQSD BAU RPS Qualifying Source Definitions'
$ws1.Range('A1').Font.Bold = $true

$ws1.Range('A2').Value2 = 'RQSD RPS Qualifying Source Definitions'
$ws1.Range('A2').Font.Bold = $true

$ws1.Range('A4').Value2 = 'Source:'
$ws1.Range('A4').Font.Bold = $true

$ws1.Range('B4').Value2 = 'see notes'
$ws1.Range('B4').Font.Bold = $false
$ws1.Range('B4').Interior.Pattern = -4142

$ws1.Range('B5').Style = 'Hyperlink'

$ws1.Range('A8').Value2 = 'Notes'
$ws1.Range('A8').Font.Bold = $true

$ws1.Range('A9').Value2 = 'Each U.S. state that has an RPS defines the sources that qualify for that RPS, leading to'
$ws1.Range('A10').Value2 = 'differences between states.  Here, we use a "clean energy standard"'
$ws1.Range('A11').Value2 = '(counting everything except fossil fuels) as our definition for the BAU case.'

$ws1.Range('A13').Value2 = 'The non-BAU version of this variable supports a boolean policy lever and is intended to be set by the'
$ws1.Range('A14').Value2 = 'model user.  The example we include uses only wind, solar, and geothermal.'
$ws1.Range('A15').Value2 = 'Hydro is excluded because of limited potential for new large hydro and land use impacts.'
$ws1.Range('A16').Value2 = 'Biomass is excluded because it is not truly carbon-neutral, and it has other issues, such as'
$ws1.Range('A17').Value2 = 'local air quality impacts and land use challenges.'
$ws1.Range('A18').Value2 = 'Nuclear is excluded because of the need to manage nuclear waste.'

$ws1.Columns.Item(2).ColumnWidth = 83.76432291666667
$ws1.Range('A1').Select()

# ---------------------------------------------------------------------------
# Shared layout for the two data sheets ("RQSD-BRQSD" and "RQSD-RQSD")
# ---------------------------------------------------------------------------
function Set-RqsdSheet($ws, $hydro, $biomass, $nuclear) {
    $ws.UsedRange.Clear()

    $ws.Range('A1').Value2 = 'Electricity Source'
    $ws.Range('A1').Font.Bold = $true
    $ws.Range('B1').Value2 = 'Qualifies for RPS (Boolean)'
    $ws.Range('B1').Font.Bold = $true
    $ws.Range('B1').HorizontalAlignment = -4152

    $ws.Range('A2').Value2 = 'hard coal'
    $ws.Range('B2').Value2 = 0

    $ws.Range('A3').Value2 = 'natural gas nonpeaker'
    $ws.Range('B3').Value2 = 0

    $ws.Range('A4').Value2 = 'nuclear'
    $ws.Range('B4').Value2 = $nuclear

    $ws.Range('A5').Value2 = 'hydro'
    $ws.Range('B5').Value2 = $hydro

    $ws.Range('A6').Value2 = 'onshore wind'
    $ws.Range('B6').Value2 = 1

    $ws.Range('A7').Value2 = 'solar PV'
    $ws.Range('B7').Value2 = 1

    $ws.Range('A8').Value2 = 'solar thermal'
    $ws.Range('B8').Value2 = 1

    $ws.Range('A9').Value2 = 'biomass'
    $ws.Range('B9').Value2 = $biomass

    $ws.Range('A10').Value2 = 'geothermal'
    $ws.Range('B10').Value2 = 1

    $ws.Range('A11').Value2 = 'petroleum'
    $ws.Range('B11').Value2 = 0

    $ws.Range('A12').Value2 = 'natural gas peaker'
    $ws.Range('B12').Value2 = 0

    $ws.Range('A13').Value2 = 'lignite'
    $ws.Range('B13').Formula = '=B2'

    $ws.Range('A14').Value2 = 'offshore wind'
    $ws.Range('B14').Value2 = 1

    $ws.Range('A15').Value2 = 'crude oil'
    $ws.Range('B15').Value2 = 0

    $ws.Range('A16').Value2 = 'heavy or residual fuel oil'
    $ws.Range('B16').Value2 = 0

    $ws.Range('A17').Value2 = 'municipal solid waste'
    $ws.Range('B17').Value2 = 0

    $ws.Columns.Item(1).ColumnWidth = 23.764322916666668
    $ws.Columns.Item(2).ColumnWidth = 22.029947916666668

    $ws.Range('B2').Select()
}

# ---------------------------------------------------------------------------
# Sheet 2: "RQSD-BRQSD" (BAU version)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item('RQSD-BRQSD')
Set-RqsdSheet $ws2 1 1 1

# ---------------------------------------------------------------------------
# Sheet 3: "RQSD-RQSD"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item('RQSD-RQSD')
Set-RqsdSheet $ws3 0 0 0

$ws1.Activate()
